$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add D12 = 0.5 (new data point for "All Primitive Data Types" row)
$ws.Range("D12").Value = 0.5

# D23: change 1 -> 2 ("Struct / Class" row)
$ws.Range("D23").Value = 2

# Add D25 = 1 ("Operator Overloading" row)
$ws.Range("D25").Value = 1

# Update the active selection to D13
$ws.Activate()
$ws.Range("D13").Select()
